# Rewrites the "colores" sheet: columns B-E get a new four-column
# layout (Color / Frase en Inglés / Traduccion al Espanol / Explicacion)
# replacing the old (Color-es / Traduccion / Color-en / English Phrase) layout.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("colores")

$data = @(
    ,@('Color', 'Frase en Inglés', 'Traducción al Español', 'Explicación')
    ,@('Aqua/Cyan', 'The sky is a clear aqua.', 'El cielo es de un azul claro.', 'Aqua es una forma más poética de decir "azul claro".')
    ,@('Aquamarine', 'The ocean is a beautiful aquamarine.', 'El océano es de un hermoso color aguamarina.', 'Aquamarine es un tono de azul verdoso, similar al color del mar.')
    ,@('Azure', 'The azure sky is so clear.', 'El cielo azur es tan claro.', 'Azure es otro sinónimo de "azul cielo", pero con una connotación más poética.')
    ,@('Beige', 'The sand is soft and beige.', 'La arena es suave y beige.', 'Beige es un color entre el crema y el marrón claro.')
    ,@('Black', 'The night sky is dark black.', 'El cielo nocturno es negro oscuro.', 'Una forma sencilla de describir el color del cielo en la noche.')
    ,@('Blue', 'The ocean is deep blue.', 'El océano es azul profundo.', 'Deep blue enfatiza la intensidad del color azul del océano.')
    ,@('Brown', 'The tree trunk is brown.', 'El tronco del árbol es marrón.', 'Un color común para los troncos de los árboles.')
    ,@('Chocolate', 'The cake is dark chocolate.', 'El pastel es de chocolate oscuro.', 'Un sabor y color común para los pasteles.')
    ,@('Coral', 'The coral reef is colorful.', 'El arrecife de coral es colorido.', 'Los arrecifes de coral son conocidos por su gran variedad de colores.')
    ,@('Crimson', 'The rose is a deep crimson.', 'La rosa es de un rojo carmesí intenso.', 'Crimson es un tono de rojo muy intenso y oscuro.')
    ,@('Dark Gray', 'The clouds are dark gray.', 'Las nubes son de un gris oscuro.', 'Un color común para las nubes antes de una tormenta.')
    ,@('Gold', 'The sun is shining gold.', 'El sol brilla dorado.', 'Una descripción poética del color del sol.')
    ,@('Gray', 'The rock is gray and old.', 'La roca es gris y vieja.', 'Un color común para las rocas y las piedras.')
    ,@('Green', 'The grass is so green.', 'La hierba está muy verde.', 'Un color asociado con la naturaleza y la vida.')
    ,@('Indigo', 'The flower is a bright indigo.', 'La flor es de un color índigo brillante.', 'Indigo es un tono de azul oscuro, casi morado.')
    ,@('Ivory', 'The dress is a soft ivory.', 'El vestido es de un suave color marfil.', 'Ivory es un color blanco cremoso, a menudo asociado con la elegancia.')
    ,@('Lavender', 'The field of lavender is purple.', 'El campo de lavanda es morado.', 'Lavender es un tono de morado pálido y suave.')
    ,@('Light Gray', 'The clouds are light gray.', 'Las nubes son de un gris claro.', 'Un color común para las nubes en un día nublado.')
    ,@('Magenta', 'The flower is a bright magenta.', 'La flor es de un color magenta brillante.', 'Magenta es un tono de rosa muy intenso y vibrante.')
    ,@('Maroon', 'The carpet is a deep maroon.', 'La alfombra es de un color burdeos intenso.', 'Maroon es un tono de rojo oscuro, casi marrón.')
    ,@('Misty Rose', 'The sky is a misty rose.', 'El cielo es de un rosa pálido.', 'Misty rose es un tono de rosa muy suave y delicado.')
    ,@('Navy', 'The sailor''s uniform is navy blue.', 'El uniforme de marinero es azul marino.', 'Navy blue es un tono de azul oscuro, asociado con la marina.')
    ,@('Olive', 'The olive oil is green.', 'El aceite de oliva es verde.', 'El aceite de oliva virgen extra tiene un color verde.')
    ,@('Orange', 'The orange is so juicy.', 'La naranja está muy jugosa.', 'Un color y sabor asociados con la fruta.')
    ,@('Pink', 'The cotton candy is pink.', 'El algodón de azúcar es rosa.', 'Un color asociado con la dulzura y la infancia.')
    ,@('Plum', 'The plum is dark purple.', 'La ciruela es de color morado oscuro.', 'Un color asociado con la fruta.')
    ,@('Purple', 'The grape is dark purple.', 'La uva es de color morado oscuro.', 'Un color común para las uvas.')
    ,@('Red', 'The apple is bright red.', 'La manzana es de un rojo brillante.', 'Un color asociado con la fruta y la energía.')
    ,@('Salmon', 'The fish is pink salmon.', 'El pescado es salmón rosado.', 'Un tipo de pescado conocido por su color rosado.')
    ,@('Silver', 'The moon is shining silver.', 'La luna brilla plateada.', 'Una descripción poética del color de la luna.')
    ,@('Sky', 'The sky is a clear blue.', 'El cielo es de un azul claro.', 'Una forma sencilla de describir el color del cielo.')
    ,@('Tan', 'The leather is a light brown.', 'El cuero es de un marrón claro.', 'Un color común para el cuero.')
    ,@('Teal', 'The duck is a bright teal.', 'El pato es de un color azul verdoso brillante.', 'Teal es un tono de azul verdoso, similar al color del agua.')
    ,@('Turquoise', 'The water is a bright turquoise.', 'El agua es de un color turquesa brillante.', 'Turquoise es un tono de azul verdoso, asociado con el agua.')
    ,@('White', 'The snow is pure white.', 'La nieve es blanca pura.', 'Un color asociado con la pureza y el invierno.')
    ,@('White Smoke', 'The smoke is white and wispy.', 'El humo es blanco y tenue.', 'Una descripción del humo blanco.')
    ,@('Yellow', 'The lemon is bright yellow.', 'El limón es de un amarillo brillante.', 'Un color asociado con la fruta y el sol.')
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 2).Value = $data[$i][0]
    $ws.Cells.Item($row, 3).Value = $data[$i][1]
    $ws.Cells.Item($row, 4).Value = $data[$i][2]
    $ws.Cells.Item($row, 5).Value = $data[$i][3]
}

# Column C widens (was merged with B at 30.57) and D widens a lot (was 11.43).
# Column B keeps its original width untouched.
$ws.Columns.Item(3).ColumnWidth = 32.666666666666664
$ws.Columns.Item(4).ColumnWidth = 39.166666666666664

# Selection moved from E2:F38 to B40 (and the view scrolled further down).
$ws.Activate()
$ws.Range("B40").Select()
